# "put instructions in the downloaded files"
#
# The workbook used to have a single data sheet ("Sheet1" holding the
# rainfall timeseries). This edit:
#   1. Adds a new first sheet named "Instructions" that holds a text-box
#      shape with an application usage guide.
#   2. Renames the original data sheet to "rain_data" and leaves its data
#      untouched.
#   3. Leaves the rain_data tab as the active/selected sheet, matching the
#      edited workbook (activeTab = the 2nd tab).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the "Instructions" sheet in front of the existing sheet ---
$instructions = $wb.Worksheets.Add()
$instructions.Name = "Instructions"

# --- 2. Rename the pre-existing data sheet (now 2nd) to "rain_data" ---
$rainData = $wb.Worksheets.Item(2)
$rainData.Name = "rain_data"

# --- 3. Add the usage-guide text box onto the Instructions sheet ---
# AddTextbox(Orientation, Left, Top, Width, Height) - Left/Top/Width/Height
# are in points; target box sits at (0pt, 1.5pt) and is 648pt x 769.5pt
# (≈ 8229599 x 9772650 EMU), matching the original authored drawing.
$tb = $instructions.Shapes.AddTextbox(1, 0, 1.5, 648, 769.5)
$tb.Name = "TextBox 1"

$guideText = @"
Application Usage Guide
=======================

This guide provides instructions for using the application and outlines the data requirements for successful analysis.

Analysis Types
--------------
IMPORTANT: The column names and tab names in the templates (both rainfall and flow) must remain unchanged for accurate analysis and app's functionality.

The application supports two types of analysis: Rainfall Analysis and Flow Analysis. In the zip package that 

Rainfall Analysis
-----------------

To perform a Rainfall Analysis, follow these steps:

1. Open the rainfall template (rainfall_template.xlsx).
2. Copy and paste the rainfall data into the template.
3. Ensure that the rain gauge values are entered into the designated rain column within the template.


Demo data is available for flow analysis under demo_data folder.

Flow Analysis
-------------

To perform a Flow Analysis, follow these steps:

1. Open the flow template (flow_template.xlsx).
2. Copy and paste the flow data from your datasheet into the downloaded template.
3. Accommodate up to four flows for analysis: inflow 1, inflow 2, outflow, and bypass.
4. Refer to the Methods tab in the template for an illustration of the possible flow type configurations.
5. You are not required to submit data for all four flow types; any combination of flow types is acceptable.

IMPORTANT: If any of the flow data types are not applicable, leave the corresponding tab blank. Do not modify the columns or the tab names.

Demo data is available for flow analysis under demo_data folder.


Data Requirements
-----------------

To ensure successful analysis, the uploaded Excel spreadsheet must conform to the following requirements:

- Flow must be reported in units of L/s, gpm, or cfs.
- The timestamp should be in 24-hour format (mm/dd/yy hh:mm:ss).
- Each tab must contain exactly two columns: one for the sample timestamps data and one for the associated values.
- The column names and the tab names must not be changed from the template.

Please ensure that your data meets these requirements before using the application for analysis.

For further assistance or inquiries, please contact our support team at stormwater@sccwrp.org
"@

$tb.TextFrame.Characters().Text = $guideText

# --- 4. Match the saved selection / active-tab state ---
# The authored file has the Instructions sheet's own cursor parked at N29
# (left over from editing the text box) while rain_data is the active tab.
$instructions.Range("N29").Select()
$rainData.Activate()
